# fixing some xlsx and cleaning nvessels tables
#
# Table5.xlsx (CDFW fish bulletin fb135) had an OCR/typo'd column header
# ("40 to G4 feet") and a stray trailing blank row left over from cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the garbled "40 to G4 feet" header in D1 -> "40 to 64 feet"
$ws.Range("D1").Value = "40 to 64 feet"

# Drop the trailing empty row (row 4) left at the bottom of the table
$ws.Rows.Item(4).Delete()

# Leave the selection on D2, matching where the cleanup left off
$ws.Range("D2").Select() | Out-Null
